# Case_1_66 (380 kV case): update res_bus/vm_pu.xlsx bus-voltage results
# for rows 2-25 (bus indices 0-23). Slack-bus voltage (col B) drops from
# 1.05 pu to 1.02 pu and every other bus voltage is recomputed accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B=1.02; C=1.010609482809288; D=1.032759647145525; E=1.012894841371706; F=1.021566736237584; I=1.032024198915173; J=1.015862738610064; K=1.035563950726745; L=1.015757626064116; M=1.024403686839626; N=1.017305380219976 }
    3 = @{ B=1.02; C=1.01195608930404; D=1.033355579263744; E=1.014044663794074; F=1.023560347132546; I=1.03226673126711; J=1.016840118094799; K=1.035969291396394; L=1.016711097299833; M=1.026200532358902; N=1.018284147695686 }
    4 = @{ B=1.02; C=1.012822566151268; D=1.033735319735869; E=1.01478488247118; F=1.024827793057758; I=1.032416675655519; J=1.017467331238495; K=1.036225025101365; L=1.017323709116644; M=1.027340343374573; N=1.018912251553988 }
    5 = @{ B=1.02; C=1.013185682875214; D=1.033893562908545; E=1.015095174267709; F=1.025355275314863; I=1.032478043611925; J=1.017729775344011; K=1.036330972626921; L=1.017580221605214; M=1.027814088870227; N=1.019175068360237 }
    6 = @{ B=1.02; C=1.013246584755431; D=1.033920050686757; E=1.015147221414558; F=1.025443529391489; I=1.032488249840295; J=1.017773768721878; K=1.036348670163254; L=1.017623231129638; M=1.027893315564966; N=1.019219124213746 }
    7 = @{ B=1.02; C=1.012827422633695; D=1.033737439682197; E=1.014789032110513; F=1.024834862268901; I=1.032417502207494; J=1.017470842870511; K=1.036226446912124; L=1.017327140679221; M=1.027346694868025; N=1.018915768172924 }
    8 = @{ B=1.02; C=1.01106559198892; D=1.032962263259591; E=1.013284222071333; F=1.022245185610508; I=1.032107614316608; J=1.016194136526184; K=1.035702296594543; L=1.016080763393022; M=1.025015701295021; N=1.017637248759147 }
    9 = @{ B=1.02; C=1.007923027341341; D=1.031551111661283; E=1.010602935523861; F=1.017506848660415; I=1.031507768688152; J=1.013903865384551; K=1.034728270913417; L=1.01385065020971; M=1.020730957623866; N=1.015343725169797 }
    10 = @{ B=1.02; C=1.005801391960342; D=1.030579556704191; E=1.008794662183051; F=1.014226939342651; I=1.031071320330002; J=1.012348870596195; K=1.034044634333437; L=1.012340363925048; M=1.017752177859602; N=1.01378652211048 }
    11 = @{ B=1.02; C=1.004876145895335; D=1.030151456997499; E=1.008006540473596; F=1.012777193719117; I=1.030873564004273; J=1.011668654375489; K=1.033740378356078; L=1.01168062022094; M=1.016432578441876; N=1.013105339904719 }
    12 = @{ B=1.02; C=1.004531459895489; D=1.029991318556566; E=1.007713008752373; F=1.012234184930459; I=1.030798781170333; J=1.011414937118939; K=1.033626116938314; L=1.011434676609137; M=1.015937879647752; N=1.012851262340563 }
    13 = @{ B=1.02; C=1.004605442213128; D=1.030025719787288; E=1.00777600826795; F=1.012350867402869; I=1.030814882567407; J=1.011469408372913; K=1.033650683002581; L=1.011487472698862; M=1.01604420083844; N=1.012905810949965 }
    14 = @{ B=1.02; C=1.004847674746568; D=1.030138242875726; E=1.007982293235923; F=1.012732400856362; I=1.030867409566748; J=1.011647703656094; K=1.03373095897295; L=1.011660308609284; M=1.016391779525666; N=1.013084359432899 }
    15 = @{ B=1.02; C=1.004996787898932; D=1.030207422928175; E=1.008109287216754; F=1.012966876510091; I=1.030899597003333; J=1.011757416926929; K=1.033780254078845; L=1.011766680666715; M=1.016605330415651; N=1.013194228509163 }
    16 = @{ B=1.02; C=1.005862657320845; D=1.03060781132487; E=1.008846857687079; F=1.014322525999694; I=1.031084259166297; J=1.012393867468966; K=1.034064652447541; L=1.012384025728117; M=1.017839122380088; N=1.013831582883971 }
    17 = @{ B=1.02; C=1.006404021554044; D=1.030856973773973; E=1.0093081310143; F=1.015164934363902; I=1.031197738018906; J=1.01279123739909; K=1.034240836167057; L=1.012769710745659; M=1.018605030454117; N=1.014229517124976 }
    18 = @{ B=1.02; C=1.006719158380726; D=1.031001591767191; E=1.009576690471909; F=1.015653453487526; I=1.031263082694505; J=1.013022352794763; K=1.034342807112825; L=1.012994117822364; M=1.019048904303736; N=1.014460960731019 }
    19 = @{ B=1.02; C=1.00682650539416; D=1.031050781903765; E=1.009668179134529; F=1.015819545747483; I=1.031285220385799; J=1.013101045148912; K=1.034377442169811; L=1.013070540954167; M=1.019199769476443; N=1.014539764837338 }
    20 = @{ B=1.02; C=1.006346003800704; D=1.030830314948629; E=1.00925869192086; F=1.015074846545949; I=1.031185650345297; J=1.012748672157461; K=1.034222015520485; L=1.012728388072897; M=1.018523152836888; N=1.014186891435821 }
    21 = @{ B=1.02; C=1.004776371323326; D=1.030105138718046; E=1.007921569360748; F=1.012620173831517; I=1.030851978395455; J=1.011595229410085; K=1.033707354210825; L=1.011609437324394; M=1.01628955218115; N=1.013031810667441 }
    22 = @{ B=1.02; C=1.003783635026667; D=1.029642688031141; E=1.007076297477288; F=1.011050698802223; I=1.030634500851759; J=1.010863902157069; K=1.033376544976065; L=1.010900776776225; M=1.014858890603523; N=1.012299444845842 }
    23 = @{ B=1.02; C=1.004310465128878; D=1.029888461844596; E=1.00752483144088; F=1.011885208917343; I=1.030750521603882; J=1.011252178499675; K=1.033552601128641; L=1.011276943476352; M=1.015619828864689; N=1.01268827258539 }
    24 = @{ B=1.02; C=1.006372221477527; D=1.030842363131121; E=1.009281032843367; F=1.015115562135691; I=1.031191114857329; J=1.012767907608284; K=1.034230522213858; L=1.012747061728903; M=1.018560158664585; N=1.01420615420319 }
    25 = @{ B=1.02; C=1.008740062934757; D=1.03192132259517; E=1.011299705771817; F=1.018752894347207; I=1.031669254552062; J=1.014500848984204; K=1.03498609155984; L=1.014431277043704; M=1.021859979569255; N=1.015941556554635 }
}

foreach ($row in ($data.Keys | Sort-Object { [int]$_ })) {
    $rowData = $data[$row]
    foreach ($col in ($rowData.Keys | Sort-Object)) {
        $colIndex = ([int][char]$col) - ([int][char]'A') + 1
        $ws.Cells.Item([int]$row, $colIndex).Value = $rowData[$col]
    }
}

Write-Output "Updated vm_pu values for $($data.Count) rows"